$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the current-year / prior-year BMV column headers
$ws.Range("G1").Value = "Curr. Year BMV"
$ws.Range("H1").Value = "Prior Year BMV"
